$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1931034482758621
$ws.Range("C2").Value = 0.5344827586206896
$ws.Range("J2").Value = 0.01724137931034483
$ws.Range("P2").Value = 0.1482758620689655
$ws.Range("S2").Value = 0.1068965517241379
$ws.Range("B3").Value = 0.01234567901234568
$ws.Range("C3").Value = 0.01851851851851852
$ws.Range("J3").Value = 0.02469135802469136
$ws.Range("P3").Value = 0.7222222222222222
$ws.Range("S3").Value = 0.2222222222222222
$ws.Range("J4").Value = 0.04651162790697674
$ws.Range("P4").Value = 0.6976744186046512
$ws.Range("S4").Value = 0.2558139534883721
$ws.Range("B6").Value = 0.06282722513089005
$ws.Range("D6").Value = 0.01047120418848168
$ws.Range("F6").Value = 0.08900523560209424
$ws.Range("J6").Value = 0.2460732984293194
$ws.Range("O6").Value = 0.02094240837696335
$ws.Range("Q6").Value = 0.1361256544502618
$ws.Range("R6").Value = 0.08376963350785341
$ws.Range("S6").Value = 0.3507853403141361
$ws.Range("B7").Value = 0.1098901098901099
$ws.Range("D7").Value = 0.01098901098901099
$ws.Range("F7").Value = 0.06593406593406594
$ws.Range("J7").Value = 0.07692307692307693
$ws.Range("O7").Value = 0.01648351648351648
$ws.Range("Q7").Value = 0.1758241758241758
$ws.Range("R7").Value = 0.08241758241758242
$ws.Range("S7").Value = 0.4615384615384616
$ws.Range("B8").Value = 0.07775377969762419
$ws.Range("D8").Value = 0.02591792656587473
$ws.Range("F8").Value = 0.06047516198704104
$ws.Range("J8").Value = 0.1209503239740821
$ws.Range("O8").Value = 0.01079913606911447
$ws.Range("Q8").Value = 0.1490280777537797
$ws.Range("R8").Value = 0.1036717062634989
$ws.Range("S8").Value = 0.4514038876889849
$ws.Range("B9").Value = 0.09722222222222222
$ws.Range("D9").Value = 0.02083333333333333
$ws.Range("F9").Value = 0.04166666666666666
$ws.Range("J9").Value = 0.0625
$ws.Range("O9").Value = 0.03472222222222222
$ws.Range("Q9").Value = 0.2430555555555556
$ws.Range("R9").Value = 0.1041666666666667
$ws.Range("S9").Value = 0.3958333333333333
$ws.Range("B10").Value = 0.1252173913043478
$ws.Range("D10").Value = 0.02086956521739131
$ws.Range("E10").Value = 0.0008695652173913044
$ws.Range("F10").Value = 0.06347826086956522
$ws.Range("J10").Value = 0.1252173913043478
$ws.Range("O10").Value = 0.01913043478260869
$ws.Range("Q10").Value = 0.1973913043478261
$ws.Range("R10").Value = 0.08695652173913043
$ws.Range("S10").Value = 0.3608695652173913
$ws.Range("G11").Value = 0.1351351351351351
$ws.Range("J11").Value = 0.07722007722007722
$ws.Range("K11").Value = 0.2007722007722008
$ws.Range("L11").Value = 0.5752895752895753
$ws.Range("S11").Value = 0.01158301158301158
$ws.Range("G12").Value = 0.7777777777777778
$ws.Range("J12").Value = 0.1633986928104575
$ws.Range("K12").Value = 0.006535947712418301
$ws.Range("L12").Value = 0.0196078431372549
$ws.Range("S12").Value = 0.03267973856209151
$ws.Range("G13").Value = 0.7837837837837838
$ws.Range("J13").Value = 0.2162162162162162
$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("J14").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.04411764705882353
$ws.Range("H15").Value = 0.1568627450980392
$ws.Range("I15").Value = 0.06372549019607843
$ws.Range("J15").Value = 0.3333333333333333
$ws.Range("K15").Value = 0.04411764705882353
$ws.Range("M15").Value = 0.01470588235294118
$ws.Range("N15").Value = 0.004901960784313725
$ws.Range("O15").Value = 0.04901960784313725
$ws.Range("S15").Value = 0.2892156862745098
$ws.Range("F16").Value = 0.01621621621621622
$ws.Range("H16").Value = 0.1945945945945946
$ws.Range("I16").Value = 0.07567567567567568
$ws.Range("J16").Value = 0.4108108108108108
$ws.Range("K16").Value = 0.07567567567567568
$ws.Range("M16").Value = 0.02702702702702703
$ws.Range("O16").Value = 0.04324324324324325
$ws.Range("S16").Value = 0.1567567567567568
$ws.Range("F17").Value = 0.005208333333333333
$ws.Range("H17").Value = 0.2109375
$ws.Range("I17").Value = 0.07291666666666667
$ws.Range("J17").Value = 0.4270833333333333
$ws.Range("K17").Value = 0.078125
$ws.Range("M17").Value = 0.01822916666666667
$ws.Range("N17").Value = 0.002604166666666667
$ws.Range("O17").Value = 0.05989583333333334
$ws.Range("S17").Value = 0.125
$ws.Range("F18").Value = 0.01030927835051546
$ws.Range("H18").Value = 0.2061855670103093
$ws.Range("I18").Value = 0.06185567010309279
$ws.Range("J18").Value = 0.4329896907216495
$ws.Range("K18").Value = 0.08762886597938144
$ws.Range("M18").Value = 0.02061855670103093
$ws.Range("O18").Value = 0.08762886597938144
$ws.Range("S18").Value = 0.09278350515463918
$ws.Range("F19").Value = 0.01410788381742739
$ws.Range("H19").Value = 0.2298755186721992
$ws.Range("I19").Value = 0.06639004149377593
$ws.Range("J19").Value = 0.3634854771784232
$ws.Range("K19").Value = 0.1112033195020747
$ws.Range("M19").Value = 0.01825726141078838
$ws.Range("N19").Value = 0.0008298755186721991
$ws.Range("O19").Value = 0.06721991701244813
$ws.Range("S19").Value = 0.1286307053941909
